$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '30.167.48'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -0.73%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.906.18'
$ws.Range('D3').Style = 'Normal'
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.002'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  +0.24%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '0.7276'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -5.53%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '242.82'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -2.03%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '1.003'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +0.31%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3116'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -3.47%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '26.41'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -6.06%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.06882'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -3.36%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.7744'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -1.64%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.07953'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -0.93%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.891.33'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -2.43%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.243'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -2.72%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '91.21'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -4.16%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '30.109.06'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -0.94%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '14.08'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -3.50%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '5.823'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -0.45%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.000007746'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -3.44%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '238.42'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -6.91%  '
$ws.Range('E21').Value = '  +0.37%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '2.152.42'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -1.95%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '1.004'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +0.47%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '6.906'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +2.03%  '
$ws.Range('E25').Value = '  -2.98%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '164.63'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +0.44%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '19.02'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -0.71%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.1269'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -5.83%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.055'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -10.90%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.352'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -1.04%  '
$ws.Range('E31').Value = '  +0.97%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.276'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -3.79%  '
$ws.Range('E33').Value = '  -2.22%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.05125'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -1.53%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.285'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -0.01%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.7355'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -2.21%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.754'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -0.73%  '
$ws.Range('E38').Value = '  -2.48%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.782'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -0.99%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '6.339'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -2.98%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '74.50'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -5.85%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.4410'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -2.83%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.930'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -2.87%  '
$ws.Range('E44').Value = '  +0.14%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.8349'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -0.31%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '100.85'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -0.71%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '7.544'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +0.30%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '9.727'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -1.37%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '37.39'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -0.29%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '2.047.60'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -2.18%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '941.52'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -4.45%  '
